$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 649, shifting existing rows 649:690 down to 650:691
$ws.Rows.Item(649).Insert()

# The newly inserted row inherits a blank format; make sure the date text in
# column A is written as plain text (not auto-converted to a date serial)
$ws.Cells.Item(649, 1).NumberFormat = "@"
$ws.Cells.Item(649, 1).Value = "2026/01/14"
$ws.Cells.Item(649, 2).Value = "水"
$ws.Cells.Item(649, 3).Value = 17
$ws.Cells.Item(649, 4).Value = 201

# Remove the formatting we applied above so the new row matches the plain,
# unstyled look of the other data rows
$ws.Rows.Item(649).ClearFormats()
